$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 345
$ws.Range("I2").Value = 971
$ws.Range("J2").Value = 4277
$ws.Range("K2").Value = 28
$ws.Range("L2").Value = 1191
$ws.Range("M2").Value = 61
$ws.Range("N2").Value = 734
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 21
$ws.Range("S2").Value = 438
$ws.Range("T2").Value = 781
$ws.Range("U2").Value = 57
$ws.Range("V2").Value = 6438
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 6533
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 99
$ws.Range("AA2").Value = 49
